$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Note: "37÷4=" must be replaced (to "31÷7=") before "36÷4=" is replaced
# to "37÷4=", otherwise the newly-introduced "37÷4=" text would collide
# with the still-pending replacement target.
Replace-Exact "37÷4=" "31÷7="

Replace-Exact "63÷5=" "82÷2="
Replace-Exact "50÷7=" "81÷9="
Replace-Exact "36÷4=" "37÷4="
Replace-Exact "18÷8=" "22÷7="
Replace-Exact "24÷4=" "81÷4="
Replace-Exact "85÷5=" "63÷8="
Replace-Exact "94÷8=" "60÷2="
Replace-Exact "49÷4=" "29÷6="
Replace-Exact "52÷5=" "30÷7="
Replace-Exact "57÷3=" "58÷5="
Replace-Exact "93÷5=" "60÷7="
Replace-Exact "61÷8=" "75÷2="
Replace-Exact "10÷6=" "56÷5="
Replace-Exact "58÷2=" "15÷7="
Replace-Exact "82÷7=" "37÷8="
Replace-Exact "18÷9=" "14÷5="
Replace-Exact "39÷6=" "49÷7="
Replace-Exact "85÷6=" "71÷8="
Replace-Exact "56÷4=" "70÷6="
Replace-Exact "55÷6=" "71÷3="
Replace-Exact "84÷2=" "82÷8="
Replace-Exact "19÷9=" "90÷4="
Replace-Exact "82÷9=" "81÷4="
Replace-Exact "82÷5=" "59÷9="
